$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.301.91"
$ws.Range("E2").Value = "  -0.76%  "

$ws.Range("D3").Value = "2.520.48"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.02"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.44"
$ws.Range("E6").Value = "  -1.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  -0.97%  "

$ws.Range("D9").Value = "2.544.18"
$ws.Range("E9").Value = "  +1.32%  "

$ws.Range("E10").Value = "  -0.70%  "

$ws.Range("E11").Value = "  -0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.59"
$ws.Range("E12").Value = "  +1.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.363"
$ws.Range("E13").Value = "  +1.94%  "

$ws.Range("D14").Value = "2.994.98"
$ws.Range("E14").Value = "  +1.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.72"
$ws.Range("E15").Value = "  -3.70%  "

$ws.Range("D16").Value = "59.296.08"
$ws.Range("E16").Value = "  -0.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").Value = "2.532.12"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.19"
$ws.Range("E19").Value = "  -1.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.28"
$ws.Range("E20").Value = "  -3.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.23"
$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("E22").Value = "  +1.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.85"
$ws.Range("E23").Value = "  -0.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.00"
$ws.Range("E24").Value = "  +1.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.434"
$ws.Range("E25").Value = "  -5.47%  "

$ws.Range("E26").Value = "  +1.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.01"
$ws.Range("E28").Value = "  +2.32%  "

$ws.Range("D29").Value = "0.0₃0785"
$ws.Range("E29").Value = "  -2.67%  "

$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.67"
$ws.Range("E31").Value = "  -2.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  -8.49%  "

$ws.Range("E33").Value = "  +3.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.68"
$ws.Range("E35").Value = "  +0.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.73"
$ws.Range("E36").Value = "  -0.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.39"
$ws.Range("E37").Value = "  -3.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.60"
$ws.Range("E38").Value = "  -9.00%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.98"
$ws.Range("E39").Value = "  +0.84%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.59"
$ws.Range("E40").Value = "  -5.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.830"
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.70"
$ws.Range("E42").Value = "  -2.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "292.11"
$ws.Range("E43").Value = "  -9.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.602"
$ws.Range("E45").Value = "  -0.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.81"
$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("E47").Value = "  -0.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.83"
$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.18"
$ws.Range("E49").Value = "  -3.70%  "

$ws.Range("E50").Value = "  -3.96%  "

$ws.Range("E51").Value = "  -2.80%  "
